# Накладная pizza: remove the "Базилик" line item (row 21) from the order
# table, fix up the two quantities that were also revised, and update the
# "amount in words" footer text to match the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Delete the entire row holding "Базилик" (item #3, 0.1 kg x 3900).
#    This shifts every row below it up by one (data rows + totals + footer).
$ws.Rows.Item(21).Delete()

# 2) The "No." column (B) in the goods table is typed in by hand, not a
#    formula, and after the deletion it must simply read 1..26 again.
for ($i = 0; $i -lt 26; $i++) {
    $ws.Cells.Item(19 + $i, 2).Value = $i + 1
}

# 3) Two quantities were revised at the same time as the row removal:
#      - "микс салата" (now row 22): 4 -> 2
#      - "Петрушка"     (now row 34): 0.2 -> 0.15
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(34, 6).Value = 0.15

# 4) Update the hand-typed "amount in words" footer (row 47 after the
#    shift) so it reflects the new item count (26) and new grand total
#    (42 989 KZT), keeping the existing rich-text run formatting.
$countCell = $ws.Cells.Item(47, 2)
$countCell.Value = "Всего отпущено количество наименований (прописью) двадцать шесть"
$countChars = $countCell.Characters(52, 13)
$countChars.Font.Italic = $true
$countChars.Font.Underline = $true
$countChars.Font.Size = 8
$countChars.Font.Name = "Arial"

$sumCell = $ws.Cells.Item(47, 5)
$sumCell.Value = " на сумму (прописью), в KZT сорок две тысячи девятьсот восемьдесят девять тенге 00 тиын"
$sumWordsChars = $sumCell.Characters(29, 53)
$sumWordsChars.Font.Italic = $true
$sumWordsChars.Font.Underline = $true
$sumWordsChars.Font.Size = 8
$sumWordsChars.Font.Name = "Arial"
$tiynChars = $sumCell.Characters(82, 8)
$tiynChars.Font.Italic = $true
$tiynChars.Font.Underline = $true
$tiynChars.Font.Size = 8
$tiynChars.Font.Name = "Arial"
$tiynChars.Font.Color = 255
